$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: ara / MNA  (mirrors row 2 "eng / MNA" structurally) ---
$ws.Range("A2:E2").Copy($ws.Range("A6:E6"))
$ws.Range("A6").Value = "ara"
$ws.Range("C6").Value = "الحكم اليدوي"
$ws.Range("D6").Value = "الرفض أثناء الحكم اليدوي"

# --- Row 7: ara / CLR  (mirrors row 3 "eng / CLR" structurally) ---
$ws.Range("A3:E3").Copy($ws.Range("A7:E7"))
$ws.Range("A7").Value = "ara"
$ws.Range("C7").Value = "رفض العميل"
$ws.Range("D7").Value = "الرفض في تسجيل العميل"

# Apply wrap-text + left/top aligned style to the new description cells,
# then propagate that exact style (format-only) to the other three cells
# so they all share a single new style record (matches the single extra
# cellXfs entry the source workbook gained).
$c = $ws.Range("C6")
$c.HorizontalAlignment = -4131
$c.VerticalAlignment = -4107
$c.WrapText = $true

$c.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("C7:D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights to match the new, taller Arabic content rows
$ws.Range("A6:E6").RowHeight = 16.4
$ws.Range("A7:E7").RowHeight = 16.4

# Column widths for the new, wider description/remarks columns
$ws.Columns("C").ColumnWidth = 30.42
$ws.Columns("D").ColumnWidth = 45.3

# Update selection to match the post-edit state
$ws.Range("C6:D7").Select()
